# Nova versao avisa envio
#
# Insert a new "CPF/CNPJ" header column between column A (DESTINATARIO)
# and the old column B (WHATSAPP (DDD + DDI + NUMERO)). The remaining
# header cells (WHATSAPP, CODIGO PEDIDO, OBJETO) shift one column to the
# right, and a matching-width column is added for the new last column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

# Remember the current header text before we start overwriting cells.
$whatsapp = $ws.Range("B1").Value()
$codigo   = $ws.Range("C1").Value()
$objeto   = $ws.Range("D1").Value()

# Shift the header formatting one column to the right (copy the cell
# format only, so the shared-string values below are untouched).
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial($xlPasteFormats)

$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial($xlPasteFormats)

$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial($xlPasteFormats)

$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial($xlPasteFormats)

$ws.Application.CutCopyMode = $false

# Now move the header text itself into its new position.
$ws.Range("E1").Value = $objeto
$ws.Range("D1").Value = $codigo
$ws.Range("C1").Value = $whatsapp
$ws.Range("B1").Value = "CPF/CNPJ"

# Give the new last column (E) the same width as column D.
$ws.Columns("E").ColumnWidth = $ws.Columns("D").ColumnWidth()
